$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new header cells for the Hint/Popup columns.
$ws.Range("G1").Value = "Hint"
$ws.Range("H1").Value = "Popup"

# Match the look of the existing header row (bold font + grey fill,
# left/top aligned) by copying the format from the neighbouring
# "Output" header cell.
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the active selection to G6, as in the edited workbook.
$ws.Range("G6").Select()
